$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Title-page author line 1: "Souvik Bhattacharyya, Abhishek Kumar Verma,"
#    -> "Souvik Bhattacharyya" (proofErr-wrapped "Souvik")
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$xml2 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Souvik</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Bhattacharyya</w:t>
  </w:r>
</w:p>
"@
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 2) Title-page author line 2 (was "Shreyash Pattnayak, Aseem Bharadwaj,
#    Achin Gupta" with spacing line=480): drop the spacing override and
#    empty the paragraph out entirely.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xml3 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 3) "Almabetter, Bangalore" -> "Almabetter" (proofErr) + ", Bangalore"
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$xml5 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Almabetter</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>, Bangalore</w:t>
  </w:r>
</w:p>
"@
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------------
# 4) Abstract paragraph: split off "fooding" with proofErr wrapping.
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs.Item(11)
$xml11 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">This </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">project deals with the trends of market in the perspective of hotel booking and how it is affected by variables such as date, time of booking, size of family etc. The study comprises of extensive research </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">on hotel bookings in two sectors, namely Resorts and City hotels. The work done is aimed at identifying various aspect and features which customers would prefer in either of those aforementioned options like </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>fooding</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>. The results that were generated from this research aim at finding the best period for booking a hotel room as per the convenience of guests. The proper implementation of the analysis would help guests find their best deals in terms of hotel bookings.</w:t>
  </w:r>
</w:p>
"@
$p11.Range.InsertXML($xml11)

# ---------------------------------------------------------------------------
# 5) "Direct, etc<tab>" -> wrap "etc" with proofErr spellStart/spellEnd.
# ---------------------------------------------------------------------------
$p58 = $d.Paragraphs.Item(58)
$xml58 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="4128"/>
    </w:tabs>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Direct</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>etc</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
  </w:r>
</w:p>
"@
$p58.Range.InsertXML($xml58)

# ---------------------------------------------------------------------------
# 6) "Libraries we used are" lost its <w:lastRenderedPageBreak/> marker in
#    the target (a pagination artifact); nothing to script here, Word
#    recomputes this automatically.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 7) "Numpy" bullet -> wrap with proofErr spellStart/spellEnd.
# ---------------------------------------------------------------------------
$p71 = $d.Paragraphs.Item(71)
$xml71 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Numpy</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
$p71.Range.InsertXML($xml71)

# ---------------------------------------------------------------------------
# 8) "Matplotlib.pyplot" bullet -> wrap with proofErr spellStart/spellEnd.
# ---------------------------------------------------------------------------
$p74 = $d.Paragraphs.Item(74)
$xml74 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Matplotlib.pyplot</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
$p74.Range.InsertXML($xml74)

# ---------------------------------------------------------------------------
# 9) Matplotlib description paragraph: split off "Tkinter" and "wxPython"
#    with proofErr wrapping.
# ---------------------------------------------------------------------------
$p75 = $d.Paragraphs.Item(75)
$xml75 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:ind w:left="720"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t xml:space="preserve">Matplotlib is a plotting library for the Python programming language and its numerical mathematics extension NumPy. It provides an object-oriented API for embedding plots into applications using general-purpose GUI toolkits like </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Tkinter</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>wxPython</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>, Qt, or GTK.</w:t>
  </w:r>
</w:p>
"@
$p75.Range.InsertXML($xml75)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
